$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string into a cell while keeping it stored as pure text
# (avoids Excel auto-converting numeric-looking strings like "0.9979" or
# "0.000009182" into real numbers). We build a temporary formula that
# evaluates to the exact text, then copy/paste-special as values only,
# which freezes the cell as a literal text value without touching the
# cell's number format / style.
function Set-TextCell {
    param(
        $Sheet,
        [string]$Address,
        [string]$Text
    )
    $escaped = $Text -replace '"', '""'
    $range = $Sheet.Range($Address)
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextCell $ws "D2" '29.044.62'
$ws.Range("E2").Value = '  +0.32%  '
Set-TextCell $ws "D3" '1.830.86'
$ws.Range("E3").Value = '  +0.39%  '
Set-TextCell $ws "D4" '0.9979'
$ws.Range("E4").Value = '  +0.00%  '
Set-TextCell $ws "D5" '241.61'
$ws.Range("E5").Value = '  -0.76%  '
Set-TextCell $ws "D6" '0.6179'
$ws.Range("E6").Value = '  -1.86%  '
Set-TextCell $ws "D7" '0.9986'
$ws.Range("E7").Value = '  +0.09%  '
Set-TextCell $ws "D8" '0.07450'
$ws.Range("E8").Value = '  -0.15%  '
Set-TextCell $ws "D9" '0.2923'
$ws.Range("E9").Value = '  -0.25%  '
Set-TextCell $ws "D10" '23.03'
$ws.Range("E10").Value = '  +0.22%  '
Set-TextCell $ws "D11" '0.07650'
$ws.Range("E11").Value = '  -0.47%  '
Set-TextCell $ws "D12" '1.823.66'
$ws.Range("E12").Value = '  -0.36%  '
Set-TextCell $ws "D13" '4.997'
$ws.Range("E13").Value = '  +0.37%  '
Set-TextCell $ws "D14" '0.6732'
$ws.Range("E14").Value = '  +1.18%  '
Set-TextCell $ws "D15" '82.80'
$ws.Range("E15").Value = '  -0.07%  '
Set-TextCell $ws "D16" '0.000009182'
$ws.Range("E16").Value = '  -4.23%  '
Set-TextCell $ws "D17" '5.891'
$ws.Range("E17").Value = '  -2.42%  '
Set-TextCell $ws "D18" '28.988.80'
$ws.Range("E18").Value = '  +0.05%  '
Set-TextCell $ws "D19" '2.075.86'
$ws.Range("E19").Value = '  -0.16%  '
Set-TextCell $ws "D20" '239.84'
$ws.Range("E20").Value = '  +6.42%  '
Set-TextCell $ws "D21" '12.67'
$ws.Range("E21").Value = '  +1.13%  '
Set-TextCell $ws "D22" '0.9989'
$ws.Range("E22").Value = '  +0.27%  '
Set-TextCell $ws "D23" '7.197'
$ws.Range("E23").Value = '  +1.08%  '
Set-TextCell $ws "D24" '0.9991'
$ws.Range("E24").Value = '  +0.05%  '
Set-TextCell $ws "D25" '158.67'
$ws.Range("E25").Value = '  -0.77%  '
Set-TextCell $ws "D26" '0.1408'
$ws.Range("E26").Value = '  -0.10%  '
Set-TextCell $ws "D27" '8.492'
$ws.Range("E27").Value = '  +0.21%  '
Set-TextCell $ws "D28" '17.87'
$ws.Range("E28").Value = '  -0.04%  '
Set-TextCell $ws "D29" '1.496'
$ws.Range("E29").Value = '  -0.02%  '
Set-TextCell $ws "D30" '0.05595'
$ws.Range("E30").Value = '  +3.08%  '
Set-TextCell $ws "D31" '4.138'
$ws.Range("E31").Value = '  +0.53%  '
Set-TextCell $ws "D32" '4.111'
$ws.Range("E32").Value = '  +1.64%  '
Set-TextCell $ws "D33" '1.198'
$ws.Range("E33").Value = '  +0.10%  '
Set-TextCell $ws "D34" '1.841'
$ws.Range("E34").Value = '  -0.30%  '
Set-TextCell $ws "D35" '0.7419'
$ws.Range("E35").Value = '  +0.07%  '
Set-TextCell $ws "D36" '1.140'
$ws.Range("E36").Value = '  +0.66%  '
Set-TextCell $ws "D37" '2.652'
$ws.Range("E37").Value = '  +1.16%  '
Set-TextCell $ws "D38" '2.769'
$ws.Range("E38").Value = '  +0.93%  '
Set-TextCell $ws "D39" '0.01783'
$ws.Range("E39").Value = '  +0.55%  '
Set-TextCell $ws "D40" '1.214.17'
$ws.Range("E40").Value = '  -1.86%  '
Set-TextCell $ws "D41" '6.415'
$ws.Range("E41").Value = '  -3.39%  '
Set-TextCell $ws "D42" '0.8953'
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("E43").Value = '  -0.01%  '
Set-TextCell $ws "D44" '101.24'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("E45").Value = '  -0.18%  '
Set-TextCell $ws "D46" '65.38'
$ws.Range("E46").Value = '  +0.72%  '
Set-TextCell $ws "D47" '0.00000000121'
$ws.Range("E47").Value = '  -1.80%  '
Set-TextCell $ws "D48" '0.5083'
$ws.Range("E48").Value = '  -0.08%  '
Set-TextCell $ws "D49" '0.4059'
$ws.Range("E49").Value = '  +0.55%  '
Set-TextCell $ws "D50" '9.129'
$ws.Range("E50").Value = '  +2.16%  '
Set-TextCell $ws "D51" '0.05810'
$ws.Range("E51").Value = '  +0.47%  '

Write-Host "Updated cryptos list"
